# Rename Sheet3 -> address
$wb = $excel.ActiveWorkbook
$wsAddress = $wb.Worksheets.Item("Sheet3")
$wsAddress.Name = "address"

# --- chat sheet: update message text (reuses the existing shared-string slot) ---
$wsChat = $wb.Worksheets.Item("chat")
$wsChat.Range("C2").Value = "add more products"

# --- Fill in the "address" sheet data ---
# Header row (row 1)
$wsAddress.Range("A1").Value = "Mobile no"
$wsAddress.Range("B1").Value = "land phone"
$wsAddress.Range("C1").Value = "street1"
$wsAddress.Range("D1").Value = "street2"
$wsAddress.Range("E1").Value = "city"
$wsAddress.Range("F1").Value = "postal code"
$wsAddress.Range("G1").Value = "land mark"

# Data row (row 2) - numeric-looking strings need a leading apostrophe
# so they are stored as text (quotePrefix style), matching the target.
$wsAddress.Range("E2").Value = "Munnar"
$wsAddress.Range("G2").Value = "church"
$wsAddress.Range("C2").Value = "KEB"
$wsAddress.Range("D2").Value = "lane"
$wsAddress.Range("A2").Value = "'997287893"
$wsAddress.Range("B2").Value = "'8310342658"
$wsAddress.Range("F2").Value = "'682001"

# Columns A and B are auto-fitted to the content (matches the <cols> block
# added to the sheet - bestFit/customWidth columns).
$wsAddress.Range("A1:B2").EntireColumn.AutoFit()

# --- chat sheet selection moves to C2 (and loses tabSelected) ---
[void]$wsChat.Range("C2").Select()

# --- address sheet becomes the active/selected sheet ---
[void]$wsAddress.Activate()
[void]$wsAddress.Range("F2").Select()

Write-Host "done"
